$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 8749.666999999999
$ws.Range("I43").Value = 7749.125
$ws.Range("J43").Value = 9550.1
$ws.Range("K43").Value = 7749.125
$ws.Range("L43").Value = 9550.1
$ws.Range("M43").Value = -7680.125
$ws.Range("N43").Value = -9688.1

# Row 70
$ws.Range("H70").Value = 1514.5
$ws.Range("I70").Value = 900
$ws.Range("J70").Value = 1637.4
$ws.Range("K70").Value = 2700
$ws.Range("L70").Value = 4912.200000000001
$ws.Range("M70").Value = -2430
$ws.Range("N70").Value = -5452.200000000001

# Row 73
$ws.Range("H73").Value = 1514.5
$ws.Range("I73").Value = 900
$ws.Range("J73").Value = 1637.4
$ws.Range("K73").Value = 2700
$ws.Range("L73").Value = 4912.200000000001
$ws.Range("M73").Value = -1764
$ws.Range("N73").Value = -6784.200000000001

# Row 86
$ws.Range("H86").Value = 457285920
$ws.Range("I86").Value = 440000300
$ws.Range("K86").Value = 440000300
$ws.Range("M86").Value = -439999177

# Row 89
$ws.Range("H89").Value = 457285920
$ws.Range("I89").Value = 440000300
$ws.Range("K89").Value = 2200001500
$ws.Range("M89").Value = -2199995884

# Row 92
$ws.Range("H92").Value = 41667436
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 98
$ws.Range("H98").Value = 2022
$ws.Range("I98").Value = 2022
$ws.Range("K98").Value = 2022
$ws.Range("M98").Value = -524

# Row 112
$ws.Range("H112").Value = 73745.42999999999
$ws.Range("J112").Value = 93523.82000000001
$ws.Range("L112").Value = 280571.46
$ws.Range("N112").Value = -282787.46

# Row 122
$ws.Range("H122").Value = 2022
$ws.Range("I122").Value = 2022
$ws.Range("K122").Value = 6066
$ws.Range("M122").Value = -3616

# Row 132
$ws.Range("H132").Value = 2496.1738
$ws.Range("I132").Value = 2377.4211
$ws.Range("K132").Value = 7132.263300000001
$ws.Range("M132").Value = -4602.263300000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 22385878
$ws.Range("I32").Value = 26043004
$ws.Range("K32").Value = 26043004
$ws.Range("M32").Value = -26042717

# Row 74
$ws.Range("H74").Value = 3128.4
$ws.Range("I74").Value = 3128.4
$ws.Range("K74").Value = 3128.4
$ws.Range("M74").Value = -2254.4

# Row 77
$ws.Range("H77").Value = 3128.4
$ws.Range("I77").Value = 3128.4
$ws.Range("K77").Value = 15642
$ws.Range("M77").Value = -11274

# Row 97
$ws.Range("H97").Value = 1165.1538
$ws.Range("I97").Value = 926.5263
$ws.Range("J97").Value = 1812.8572
$ws.Range("K97").Value = 926.5263
$ws.Range("L97").Value = 1812.8572
$ws.Range("M97").Value = -430.5263
$ws.Range("N97").Value = -2804.8572

# Row 106
$ws.Range("H106").Value = 21680
$ws.Range("I106").Value = 18000
$ws.Range("J106").Value = 22906.666
$ws.Range("K106").Value = 18000
$ws.Range("L106").Value = 22906.666
$ws.Range("M106").Value = -16738
$ws.Range("N106").Value = -25430.666

# Row 110
$ws.Range("H110").Value = 2507.4
$ws.Range("I110").Value = 1263.875
$ws.Range("K110").Value = 1263.875
$ws.Range("M110").Value = 781.125

# Row 132
$ws.Range("H132").Value = 2934.0688
$ws.Range("I132").Value = 2743.92
$ws.Range("K132").Value = 8231.76
$ws.Range("M132").Value = -5701.76

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1815.569
$ws.Range("I31").Value = 1319.7307
$ws.Range("K31").Value = 1319.7307
$ws.Range("M31").Value = -1024.7307

# Row 34
$ws.Range("H34").Value = 1815.569
$ws.Range("I34").Value = 1319.7307
$ws.Range("K34").Value = 1319.7307
$ws.Range("M34").Value = -1117.7307

# Row 105
$ws.Range("H105").Value = 2292.9
$ws.Range("I105").Value = 2178.625
$ws.Range("K105").Value = 2178.625
$ws.Range("M105").Value = -431.625

# Row 132
$ws.Range("H132").Value = 3700.2778
$ws.Range("I132").Value = 3131.8928
$ws.Range("K132").Value = 9395.678400000001
$ws.Range("M132").Value = -6865.678400000001

$ws = $wb.Worksheets.Item("CUL")
# Row 31
$ws.Range("H31").Value = 1880.6
$ws.Range("I31").Value = 1801
$ws.Range("K31").Value = 5403
$ws.Range("M31").Value = -5115

# Row 68
$ws.Range("H68").Value = 1983.48
$ws.Range("I68").Value = 1489.3
$ws.Range("J68").Value = 2312.9333
$ws.Range("K68").Value = 4467.9
$ws.Range("L68").Value = 6938.7999
$ws.Range("M68").Value = -3656.9
$ws.Range("N68").Value = -8560.7999

# Row 71
$ws.Range("H71").Value = 1983.48
$ws.Range("I71").Value = 1489.3
$ws.Range("J71").Value = 2312.9333
$ws.Range("K71").Value = 13403.7
$ws.Range("L71").Value = 20816.3997
$ws.Range("M71").Value = -9347.699999999999
$ws.Range("N71").Value = -28928.3997

# Row 81
$ws.Range("H81").Value = 2628
$ws.Range("I81").Value = 1974
$ws.Range("K81").Value = 5922
$ws.Range("M81").Value = -4799

# Row 84
$ws.Range("H84").Value = 2628
$ws.Range("I84").Value = 1974
$ws.Range("K84").Value = 17766
$ws.Range("M84").Value = -12150

# Row 129
$ws.Range("H129").Value = 1977.4375
$ws.Range("J129").Value = 2075.4666
$ws.Range("L129").Value = 6226.399800000001
$ws.Range("N129").Value = -16226.3998

# Row 139
$ws.Range("H139").Value = 2864.6
$ws.Range("I139").Value = 2864.6
$ws.Range("K139").Value = 8593.799999999999
$ws.Range("M139").Value = -3453.799999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 24925.334
$ws.Range("J93").Value = 24925.334
$ws.Range("L93").Value = 24925.334
$ws.Range("N93").Value = -28669.334

# Row 98
$ws.Range("H98").Value = 116000
$ws.Range("J98").Value = 116000
$ws.Range("L98").Value = 116000
$ws.Range("M98").Value = -121990

# Row 102
$ws.Range("H102").Value = 2270.4375
$ws.Range("I102").Value = 2187.4167
$ws.Range("J102").Value = 2519.5
$ws.Range("K102").Value = 2187.4167
$ws.Range("L102").Value = 2519.5
$ws.Range("M102").Value = -565.4167000000002
$ws.Range("N102").Value = -5763.5

# Row 132
$ws.Range("H132").Value = 2784.0715
$ws.Range("I132").Value = 1998.4286
$ws.Range("J132").Value = 3569.7144
$ws.Range("K132").Value = 5995.2858
$ws.Range("L132").Value = 10709.1432
$ws.Range("M132").Value = -3465.2858
$ws.Range("N132").Value = -15769.1432

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 2749.0435
$ws.Range("I46").Value = 1310.25
$ws.Range("J46").Value = 3051.9473
$ws.Range("K46").Value = 1310.25
$ws.Range("L46").Value = 3051.9473
$ws.Range("M46").Value = -1122.25
$ws.Range("N46").Value = -3427.9473

# Row 55
$ws.Range("H55").Value = 1912.2
$ws.Range("I55").Value = 2254.1428
$ws.Range("J55").Value = 1114.3334
$ws.Range("K55").Value = 2254.1428
$ws.Range("L55").Value = 1114.3334
$ws.Range("M55").Value = -2081.1428
$ws.Range("N55").Value = -1460.3334

# Row 93
$ws.Range("H93").Value = 125001980
$ws.Range("I93").Value = 166668320
$ws.Range("K93").Value = 166668320
$ws.Range("M93").Value = -166667072

# Row 95
$ws.Range("H95").Value = 30000
$ws.Range("I95").Value = 30000
$ws.Range("K95").Value = 30000
$ws.Range("M95").Value = -27254

# Row 132
$ws.Range("H132").Value = 15068.714
$ws.Range("I132").Value = 14269.272
$ws.Range("K132").Value = 42807.81600000001
$ws.Range("M132").Value = -40277.81600000001

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 7500000
$ws.Range("I4").Value = 10000000
$ws.Range("K4").Value = 10000000
$ws.Range("M4").Value = -9999887

# Row 33
$ws.Range("H33").Value = 25353
$ws.Range("J33").Value = 25353
$ws.Range("L33").Value = 25353
$ws.Range("N33").Value = -25853

# Row 36
$ws.Range("H36").Value = 25353
$ws.Range("J36").Value = 25353
$ws.Range("L36").Value = 25353
$ws.Range("N36").Value = -25853
